$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150, shifting existing rows 150-180 down to 151-181.
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with the new record.
$ws.Cells.Item(150, 1).Value = 8
$ws.Cells.Item(150, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(150, 3).Value = "Coquimbo"
$ws.Cells.Item(150, 4).Value = "2023-09-04"
$ws.Cells.Item(150, 5).Value = 4
$ws.Cells.Item(150, 6).Value = 100114007
$ws.Cells.Item(150, 7).Value = "Jengibre"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 400
$ws.Cells.Item(150, 11).Value = 17000
$ws.Cells.Item(150, 12).Value = 18000
$ws.Cells.Item(150, 13).Value = 17500
$ws.Cells.Item(150, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(150, 15).Value = "Perú"
$ws.Cells.Item(150, 16).Value = 1346
$ws.Cells.Item(150, 17).Value = 13
$ws.Cells.Item(150, 18).Value = "Hortaliza"
